$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk update of existing Qtd. Processos (column B) values that changed without row shift
$ws.Cells.Item(7, 2).Value = 55457
$ws.Cells.Item(26, 2).Value = 177
$ws.Cells.Item(27, 2).Value = 584
$ws.Cells.Item(28, 2).Value = 447
$ws.Cells.Item(49, 2).Value = 111
$ws.Cells.Item(62, 2).Value = 12
$ws.Cells.Item(63, 2).Value = 35
$ws.Cells.Item(72, 2).Value = 282
$ws.Cells.Item(78, 2).Value = 435
$ws.Cells.Item(81, 2).Value = 290
$ws.Cells.Item(85, 2).Value = 3
$ws.Cells.Item(86, 2).Value = 31
$ws.Cells.Item(87, 2).Value = 184
$ws.Cells.Item(88, 2).Value = 56
$ws.Cells.Item(94, 2).Value = 30
$ws.Cells.Item(96, 2).Value = 322
$ws.Cells.Item(106, 2).Value = 34
$ws.Cells.Item(116, 2).Value = 4463
$ws.Cells.Item(121, 2).Value = 4796
$ws.Cells.Item(128, 2).Value = 731
$ws.Cells.Item(131, 2).Value = 3533
$ws.Cells.Item(135, 2).Value = 585
$ws.Cells.Item(139, 2).Value = 2
$ws.Cells.Item(140, 2).Value = 46
$ws.Cells.Item(141, 2).Value = 1683
$ws.Cells.Item(145, 2).Value = 2004
$ws.Cells.Item(160, 2).Value = 2903
$ws.Cells.Item(169, 2).Value = 10395
$ws.Cells.Item(182, 2).Value = 1298
$ws.Cells.Item(202, 2).Value = 1817
$ws.Cells.Item(231, 2).Value = 413
$ws.Cells.Item(234, 2).Value = 892
$ws.Cells.Item(250, 2).Value = 420198
$ws.Cells.Item(251, 2).Value = 261386
$ws.Cells.Item(256, 2).Value = 431
$ws.Cells.Item(269, 2).Value = 79746
$ws.Cells.Item(270, 2).Value = 160916
$ws.Cells.Item(273, 2).Value = 56
$ws.Cells.Item(289, 2).Value = 267
$ws.Cells.Item(290, 2).Value = 5
$ws.Cells.Item(307, 2).Value = 152
$ws.Cells.Item(309, 2).Value = 3361
$ws.Cells.Item(314, 2).Value = 94
$ws.Cells.Item(317, 2).Value = 1120
$ws.Cells.Item(322, 2).Value = 17349
$ws.Cells.Item(333, 2).Value = 20672
$ws.Cells.Item(334, 2).Value = 169
$ws.Cells.Item(336, 2).Value = 19
$ws.Cells.Item(337, 2).Value = 10
$ws.Cells.Item(352, 2).Value = 7
$ws.Cells.Item(365, 2).Value = 162
$ws.Cells.Item(370, 2).Value = 4159
$ws.Cells.Item(375, 2).Value = 1451
$ws.Cells.Item(381, 2).Value = 77
$ws.Cells.Item(412, 2).Value = 593
$ws.Cells.Item(433, 2).Value = 239
$ws.Cells.Item(463, 2).Value = 805
$ws.Cells.Item(482, 2).Value = 10
$ws.Cells.Item(499, 2).Value = 348
$ws.Cells.Item(503, 2).Value = 7
$ws.Cells.Item(504, 2).Value = 186
$ws.Cells.Item(508, 2).Value = 88
$ws.Cells.Item(509, 2).Value = 648
$ws.Cells.Item(515, 2).Value = 1130
$ws.Cells.Item(527, 2).Value = 4812
$ws.Cells.Item(541, 2).Value = 35
$ws.Cells.Item(549, 2).Value = 82
$ws.Cells.Item(551, 2).Value = 20
$ws.Cells.Item(566, 2).Value = 1751
$ws.Cells.Item(571, 2).Value = 1755
$ws.Cells.Item(581, 2).Value = 6507
$ws.Cells.Item(618, 2).Value = 16
$ws.Cells.Item(619, 2).Value = 25035
$ws.Cells.Item(624, 2).Value = 176
$ws.Cells.Item(625, 2).Value = 429
$ws.Cells.Item(632, 2).Value = 191
$ws.Cells.Item(636, 2).Value = 11
$ws.Cells.Item(640, 2).Value = 2694
$ws.Cells.Item(645, 2).Value = 14
$ws.Cells.Item(676, 2).Value = 2209
$ws.Cells.Item(679, 2).Value = 105
$ws.Cells.Item(698, 2).Value = 3233
$ws.Cells.Item(699, 2).Value = 16
$ws.Cells.Item(700, 2).Value = 39
$ws.Cells.Item(715, 2).Value = 114
$ws.Cells.Item(732, 2).Value = 2
$ws.Cells.Item(738, 2).Value = 15862
$ws.Cells.Item(758, 2).Value = 938
$ws.Cells.Item(771, 2).Value = 1120
$ws.Cells.Item(793, 2).Value = 1050
$ws.Cells.Item(803, 2).Value = 24
$ws.Cells.Item(810, 2).Value = 231
$ws.Cells.Item(814, 2).Value = 697
$ws.Cells.Item(821, 2).Value = 1126
$ws.Cells.Item(824, 2).Value = 911
$ws.Cells.Item(831, 2).Value = 7
$ws.Cells.Item(842, 2).Value = 6210
$ws.Cells.Item(846, 2).Value = 1518
$ws.Cells.Item(847, 2).Value = 204
$ws.Cells.Item(849, 2).Value = 11
$ws.Cells.Item(862, 2).Value = 4689
$ws.Cells.Item(877, 2).Value = 1466
$ws.Cells.Item(880, 2).Value = 1405
$ws.Cells.Item(889, 2).Value = 1962
$ws.Cells.Item(898, 2).Value = 335
$ws.Cells.Item(906, 2).Value = 235
$ws.Cells.Item(921, 2).Value = 19229
$ws.Cells.Item(932, 2).Value = 125
$ws.Cells.Item(938, 2).Value = 14580
$ws.Cells.Item(942, 2).Value = 33
$ws.Cells.Item(944, 2).Value = 16
$ws.Cells.Item(946, 2).Value = 9756
$ws.Cells.Item(951, 2).Value = 2296
$ws.Cells.Item(952, 2).Value = 331
$ws.Cells.Item(954, 2).Value = 4814
$ws.Cells.Item(957, 2).Value = 5563
$ws.Cells.Item(960, 2).Value = 99
$ws.Cells.Item(976, 2).Value = 90
$ws.Cells.Item(979, 2).Value = 62321
$ws.Cells.Item(996, 2).Value = 5
$ws.Cells.Item(1008, 2).Value = 6
$ws.Cells.Item(1016, 2).Value = 29
$ws.Cells.Item(1028, 2).Value = 1392
$ws.Cells.Item(1031, 2).Value = 719
$ws.Cells.Item(1032, 2).Value = 74
$ws.Cells.Item(1038, 2).Value = 31932
$ws.Cells.Item(1039, 2).Value = 4
$ws.Cells.Item(1043, 2).Value = 50344
$ws.Cells.Item(1050, 2).Value = 382
$ws.Cells.Item(1065, 2).Value = 107
$ws.Cells.Item(1069, 2).Value = 40
$ws.Cells.Item(1079, 2).Value = 320

# Row 1086 (STF RG 1179 / 2) becomes "STF RG 1178" with value 10
$ws.Cells.Item(1086, 1).Value = "STF RG 1178"
$ws.Cells.Item(1086, 2).Value = 10

# Insert new row 1087: "STF RG 1179" / 2 (the old content of row 1086, now shifted down)
$ws.Range("A1087:B1087").Insert(-4121)
$ws.Range("A1086:B1086").Copy()
$ws.Range("A1087:B1087").PasteSpecial(-4122)
$ws.Cells.Item(1087, 1).Value = "STF RG 1179"
$ws.Cells.Item(1087, 2).Value = 2
$ws.Rows.Item(1087).RowHeight = 32.517578125

# Insert 3 new rows before the final row (previously 1087, "STF RG 1189", now at 1091)
# New rows: STF RG 1184 / 92, STF RG 1186 / 8, STF RG 1187 / 1
$ws.Range("A1088:B1090").Insert(-4121)
$ws.Range("A1087:B1087").Copy()
$ws.Range("A1088:B1090").PasteSpecial(-4122)
$ws.Range("A1088:B1090").RowHeight = 32.517578125

$ws.Cells.Item(1088, 1).Value = "STF RG 1184"
$ws.Cells.Item(1088, 2).Value = 92
$ws.Cells.Item(1089, 1).Value = "STF RG 1186"
$ws.Cells.Item(1089, 2).Value = 8
$ws.Cells.Item(1090, 1).Value = "STF RG 1187"
$ws.Cells.Item(1090, 2).Value = 1

# The row that was originally 1087 ("STF RG 1189"/1) is now at row 1091 - unchanged values, no action needed

$ws.Application.CutCopyMode = $false
